$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04717506928171853
$ws.Range("C2").Value = 1.491543766192013
$ws.Range("D2").Value = 8.455526954849582
$ws.Range("E2").Value = 2.907838880483164
$ws.Range("F2").Value = 2.975876279823447
$ws.Range("G2").Value = 22
$ws.Range("B3").Value = -0.4214464237818173
$ws.Range("C3").Value = 1.016004021736437
$ws.Range("D3").Value = 3.540968988508833
$ws.Range("E3").Value = 1.881746260394539
$ws.Range("F3").Value = 1.879233752200579
$ws.Range("G3").Value = 21
$ws.Range("B4").Value = -0.04213737398666674
$ws.Range("C4").Value = 0.7820104901831507
$ws.Range("D4").Value = 1.001266313089072
$ws.Range("E4").Value = 1.000632956227743
$ws.Range("F4").Value = 1.025717081543998
$ws.Range("G4").Value = 20
$ws.Range("B5").Value = -0.0547463289693642
$ws.Range("C5").Value = 0.7907933146236926
$ws.Range("D5").Value = 1.075494849714573
$ws.Range("E5").Value = 1.037060677932865
$ws.Range("F5").Value = 1.063992901354768
$ws.Range("G5").Value = 19
$ws.Range("B6").Value = -0.1059385590168174
$ws.Range("C6").Value = 0.6529650943868757
$ws.Range("D6").Value = 0.7643346331479901
$ws.Range("E6").Value = 0.8742623365718039
$ws.Range("F6").Value = 0.8929794736955126
$ws.Range("G6").Value = 18
$ws.Range("B7").Value = -0.06115313360779569
$ws.Range("C7").Value = 0.7409233030117504
$ws.Range("D7").Value = 0.8853216378027108
$ws.Range("E7").Value = 0.9409153191455174
$ws.Range("F7").Value = 0.9678227124871316
$ws.Range("G7").Value = 17
$ws.Range("B8").Value = 0.04633068755579241
$ws.Range("C8").Value = 0.7066816895393113
$ws.Range("D8").Value = 0.7257841232605365
$ws.Range("E8").Value = 0.8519296468961135
$ws.Range("F8").Value = 0.8785670701173324
$ws.Range("G8").Value = 16
$ws.Range("B9").Value = -0.03682260890791801
$ws.Range("C9").Value = 0.5844290800941564
$ws.Range("D9").Value = 0.7081477771230911
$ws.Range("E9").Value = 0.841515167494378
$ws.Range("F9").Value = 0.8702166433441642
$ws.Range("G9").Value = 15
$ws.Range("B10").Value = -0.00294460816098668
$ws.Range("C10").Value = 0.5209809220943237
$ws.Range("D10").Value = 0.4617339020710022
$ws.Range("E10").Value = 0.6795100455997706
$ws.Range("F10").Value = 0.7051542787521982
$ws.Range("G10").Value = 14
$ws.Range("B11").Value = -0.04202418731684543
$ws.Range("C11").Value = 0.6363847651420124
$ws.Range("D11").Value = 0.7103350813464621
$ws.Range("E11").Value = 0.8428137880614329
$ws.Range("F11").Value = 0.8761372436890558
$ws.Range("G11").Value = 13

Write-Host "Updated ifoCAST full series evaluation rows."
